# Add a new trade row (row 9) to the IBB bag-trade log, mirroring the
# format of the existing rows above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 9895.27
$ws.Range("B9").Value = 9974.07
$ws.Range("C9").Value = 282.89999999999998
$ws.Range("D9").Value = 285.14
$ws.Range("E9").Value = $true
$ws.Range("F9").Value = 0.79
$ws.Range("G9").Value = 42609.487245370372
$ws.Range("H9").Value = $false

# Match the date/time number formatting used by the other rows in column G.
$ws.Range("G3").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = 42609.487245370372
